{"js": "// Update the date line and every answer cell in the addition/subtraction\n// worksheet table, in document order (the date paragraph first, then each\n// table-cell paragraph left-to-right, top-to-bottom). The document structure\n// (paragraphs/rows/cells/runs/formatting) is unchanged; only each run's text\n// content is replaced.\nconst newTexts = [\n  \"2025-05-27 Tuesday\",\n  \"40+18=58\",\n  \"93-31=62\",\n  \"99-5=94\",\n  \"77-33=44\",\n  \"88-3=85\",\n  \"48+41=89\",\n  \"28+62=90\",\n  \"62-57=5\",\n  \"47-30=17\",\n  \"77-53=24\",\n  \"29-11=18\",\n  \"85-69=16\",\n  \"27-0=27\",\n  \"81-81=0\",\n  \"60+20=80\",\n  \"88-34=54\",\n  \"43-26=17\",\n  \"93-1=92\",\n  \"64-12=52\",\n  \"10+2=12\",\n  \"56-6=50\",\n  \"64-25=39\",\n  \"61+34=95\",\n  \"71+4=75\",\n  \"38-19=19\",\n  \"46-1=45\",\n  \"53-25=28\",\n  \"43+9=52\",\n  \"67-11=56\",\n  \"31+28=59\",\n  \"34+39=73\",\n  \"11+48=59\",\n  \"77-51=26\",\n  \"70+24=94\",\n  \"71-4=67\",\n  \"24+15=39\",\n  \"21+54=75\",\n  \"52-18=34\",\n  \"1+84=85\",\n  \"16+28=44\",\n  \"56+5=61\",\n  \"23-4=19\",\n  \"19-18=1\",\n  \"19+25=44\",\n  \"10+74=84\",\n  \"58-47=11\",\n  \"40+29=69\",\n  \"61+13=74\",\n  \"60-38=22\",\n  \"25+27=52\",\n  \"15+61=76\",\n  \"23-10=13\",\n  \"75+10=85\",\n  \"47+32=79\",\n  \"42+23=65\",\n  \"27+66=93\",\n  \"83-16=67\",\n  \"81-14=67\",\n  \"43-43=0\",\n  \"60-6=54\",\n  \"4+24=28\",\n  \"45+41=86\",\n  \"83-12=71\",\n  \"59-39=20\",\n  \"34+40=74\",\n  \"3+1=4\",\n  \"31+41=72\",\n  \"0+23=23\",\n  \"32+27=59\",\n  \"41-8=33\",\n  \"80-64=16\",\n  \"32+9=41\",\n  \"90-42=48\",\n  \"65+23=88\",\n  \"56+31=87\",\n  \"13+70=83\",\n  \"84-9=75\",\n  \"0+54=54\",\n  \"3+11=14\",\n  \"10+33=43\",\n  \"76-24=52\",\n  \"22+7=29\",\n  \"50-9=41\",\n  \"50+40=90\",\n  \"76-26=50\",\n  \"45+32=77\",\n  \"50+46=96\",\n  \"67-40=27\",\n  \"52-43=9\",\n  \"49+2=51\",\n  \"67-24=43\",\n  \"52-33=19\",\n  \"36-18=18\",\n  \"85-79=6\",\n  \"20+63=83\",\n  \"31+30=61\",\n  \"0+34=34\",\n  \"39-24=15\",\n  \"28+67=95\",\n  \"66-16=50\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== newTexts.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + newTexts.length + \" got \" + items.length\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every answer cell in the addition/subtraction\n# worksheet table (20 rows x 5 columns), in document order. Document\n# structure (paragraphs/rows/cells/runs/formatting) is unchanged; only\n# each cell/paragraph's text content is replaced.\n\n$d = $word.ActiveDocument\n\n$newDate = \"2025-05-27 Tuesday\"\n$d.Paragraphs.Item(1).Range.Text = $newDate\n\n$newValues = @(\n    @(\"40+18=58\", \"93-31=62\", \"99-5=94\", \"77-33=44\", \"88-3=85\"),\n    @(\"48+41=89\", \"28+62=90\", \"62-57=5\", \"47-30=17\", \"77-53=24\"),\n    @(\"29-11=18\", \"85-69=16\", \"27-0=27\", \"81-81=0\", \"60+20=80\"),\n    @(\"88-34=54\", \"43-26=17\", \"93-1=92\", \"64-12=52\", \"10+2=12\"),\n    @(\"56-6=50\", \"64-25=39\", \"61+34=95\", \"71+4=75\", \"38-19=19\"),\n    @(\"46-1=45\", \"53-25=28\", \"43+9=52\", \"67-11=56\", \"31+28=59\"),\n    @(\"34+39=73\", \"11+48=59\", \"77-51=26\", \"70+24=94\", \"71-4=67\"),\n    @(\"24+15=39\", \"21+54=75\", \"52-18=34\", \"1+84=85\", \"16+28=44\"),\n    @(\"56+5=61\", \"23-4=19\", \"19-18=1\", \"19+25=44\", \"10+74=84\"),\n    @(\"58-47=11\", \"40+29=69\", \"61+13=74\", \"60-38=22\", \"25+27=52\"),\n    @(\"15+61=76\", \"23-10=13\", \"75+10=85\", \"47+32=79\", \"42+23=65\"),\n    @(\"27+66=93\", \"83-16=67\", \"81-14=67\", \"43-43=0\", \"60-6=54\"),\n    @(\"4+24=28\", \"45+41=86\", \"83-12=71\", \"59-39=20\", \"34+40=74\"),\n    @(\"3+1=4\", \"31+41=72\", \"0+23=23\", \"32+27=59\", \"41-8=33\"),\n    @(\"80-64=16\", \"32+9=41\", \"90-42=48\", \"65+23=88\", \"56+31=87\"),\n    @(\"13+70=83\", \"84-9=75\", \"0+54=54\", \"3+11=14\", \"10+33=43\"),\n    @(\"76-24=52\", \"22+7=29\", \"50-9=41\", \"50+40=90\", \"76-26=50\"),\n    @(\"45+32=77\", \"50+46=96\", \"67-40=27\", \"52-43=9\", \"49+2=51\"),\n    @(\"67-24=43\", \"52-33=19\", \"36-18=18\", \"85-79=6\", \"20+63=83\"),\n    @(\"31+30=61\", \"0+34=34\", \"39-24=15\", \"28+67=95\", \"66-16=50\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n\n"}
